$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "112.96") must keep
# their original text type, matching the source inlineStr cells in the diff.
# Force those specific cells to Text format before writing the value so Excel
# does not auto-convert them into numeric cells.
$textCells = @("D5","D6","D7","D9","D10","D11","D12","D15","D19","D21","D23","D24","D25","D27","D28","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D42","D44","D45","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '45.344.42'
$ws.Range('E2').Value = '  +6.17%  '
$ws.Range('D3').Value = '2.390.07'
$ws.Range('E3').Value = '  +4.87%  '
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('D5').Value = '112.96'
$ws.Range('E5').Value = '  +9.60%  '
$ws.Range('D6').Value = '318.69'
$ws.Range('E6').Value = '  +2.97%  '
$ws.Range('D7').Value = '0.637'
$ws.Range('E7').Value = '  +2.91%  '
$ws.Range('E8').Value = '  -0.37%  '
$ws.Range('D9').Value = '0.631'
$ws.Range('E9').Value = '  +5.45%  '
$ws.Range('D10').Value = '42.50'
$ws.Range('E10').Value = '  +10.76%  '
$ws.Range('D11').Value = '0.0933'
$ws.Range('E11').Value = '  +4.17%  '
$ws.Range('D12').Value = '8.71'
$ws.Range('E12').Value = '  +6.53%  '
$ws.Range('E13').Value = '  +5.68%  '
$ws.Range('E14').Value = '  +1.45%  '
$ws.Range('D15').Value = '15.86'
$ws.Range('E15').Value = '  +5.94%  '
$ws.Range('D16').Value = '2.750.40'
$ws.Range('E16').Value = '  +4.82%  '
$ws.Range('D17').Value = '2.386.08'
$ws.Range('E17').Value = '  +4.81%  '
$ws.Range('D18').Value = '45.297.45'
$ws.Range('E18').Value = '  +6.32%  '
$ws.Range('D19').Value = '7.69'
$ws.Range('E19').Value = '  +6.40%  '
$ws.Range('E20').Value = '  +4.36%  '
$ws.Range('D21').Value = '13.10'
$ws.Range('E21').Value = '  +1.39%  '
$ws.Range('E22').Value = '  +3.48%  '
$ws.Range('D23').Value = '3.55'
$ws.Range('E23').Value = '  +5.65%  '
$ws.Range('D24').Value = '269.30'
$ws.Range('E24').Value = '  +3.07%  '
$ws.Range('D25').Value = '2.33'
$ws.Range('E25').Value = '  +8.05%  '
$ws.Range('E26').Value = '  -0.60%  '
$ws.Range('D27').Value = '11.32'
$ws.Range('E27').Value = '  +7.02%  '
$ws.Range('D28').Value = '7.51'
$ws.Range('E28').Value = '  +9.81%  '
$ws.Range('E29').Value = '  +0.45%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = '39.08'
$ws.Range('E30').Value = '  +10.39%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '22.96'
$ws.Range('E31').Value = '  +4.19%  '
$ws.Range('D32').Value = '0.0949'
$ws.Range('E32').Value = '  +11.94%  '
$ws.Range('D33').Value = '169.98'
$ws.Range('E33').Value = '  +3.70%  '
$ws.Range('D34').Value = '2.98'
$ws.Range('E34').Value = '  +16.71%  '
$ws.Range('D35').Value = '0.134'
$ws.Range('E35').Value = '  +3.77%  '
$ws.Range('D36').Value = '4.95'
$ws.Range('E36').Value = '  +11.30%  '
$ws.Range('D37').Value = '0.119'
$ws.Range('E37').Value = '  +7.87%  '
$ws.Range('D38').Value = '3.13'
$ws.Range('E38').Value = '  +16.04%  '
$ws.Range('D39').Value = '0.0366'
$ws.Range('E39').Value = '  +6.05%  '
$ws.Range('E40').Value = '  +8.83%  '
$ws.Range('E41').Value = '  +13.23%  '
$ws.Range('D42').Value = '104.73'
$ws.Range('E42').Value = '  +5.68%  '
$ws.Range('E43').Value = '  +7.42%  '
$ws.Range('D44').Value = '13.64'
$ws.Range('E44').Value = '  +15.01%  '
$ws.Range('D45').Value = '71.29'
$ws.Range('E45').Value = '  +4.83%  '
$ws.Range('E46').Value = '  -0.53%  '
$ws.Range('D47').Value = '117.82'
$ws.Range('E47').Value = '  +7.72%  '
$ws.Range('D48').Value = '5.81'
$ws.Range('E48').Value = '  +13.69%  '
$ws.Range('D49').Value = '1.65'
$ws.Range('E49').Value = '  +21.23%  '
$ws.Range('D50').Value = '9.38'
$ws.Range('E50').Value = '  +9.36%  '
$ws.Range('B51').Value = 'TheGraph'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D51').Value = '0.224'
$ws.Range('E51').Value = '  +19.33%  '
